$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich text) ---
# A8: "Volume 32   Number  42" -> "...43"  (volume/number footer)
$ws.Range("A8").Value = "Volume 32   Number  43"
# C9: report date range bumped one week forward
$ws.Range("C9").Value = "Report Covering the Week  10/20/2025  Through  10/26/2025"

# --- Row 22 (72nd Pct / Shooting Vic. row) column restructuring ---
# Before: C22="0"(text) D22=2(num) E22=-100(num) F22="0"(text) G22=4 ...
# After:  C22=2(num)    D22="0"(text) E22="***.*"(text) F22=2(num) G22=3 ...
# Reassign styles by pasting formats from donor cells with the matching
# style index, then set the values (forcing text cells via NumberFormat
# "@" so "0" / "***.*" are stored as text, not as numbers).
$ws.Range("G22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 2

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"

$ws.Range("G22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2

$ws.Range("A22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 31
$ws.Range("K22").Value = 47.619047619047
$ws.Range("L22").Value = 14.814814814814
$ws.Range("M22").Value = 24

# --- Remaining weekly/28-day/YTD crime-count & %Chg figures (rows 14-30) ---

# Row 14
$ws.Range("N14").Value = -94.117647058823

# Row 15
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 26.666666666666
$ws.Range("M15").Value = 5.555555555555
$ws.Range("N15").Value = -32.142857142857

# Row 16
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -42.105263157894
$ws.Range("I16").Value = 88
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = -45
$ws.Range("L16").Value = -30.708661417322
$ws.Range("M16").Value = -54.639175257732
$ws.Range("N16").Value = -89.536266349583

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -66.666666666666
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -42.857142857142
$ws.Range("I17").Value = 249
$ws.Range("J17").Value = 277
$ws.Range("K17").Value = -10.108303249097
$ws.Range("L17").Value = 5.063291139240
$ws.Range("M17").Value = 56.603773584905
$ws.Range("N17").Value = -43.665158371040

# Row 18
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -30.769230769230
$ws.Range("I18").Value = 85
$ws.Range("J18").Value = 142
$ws.Range("K18").Value = -40.140845070422
$ws.Range("L18").Value = -22.018348623853
$ws.Range("M18").Value = -60.648148148148
$ws.Range("N18").Value = -91.033755274261

# Row 19
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -17.948717948717
$ws.Range("I19").Value = 286
$ws.Range("J19").Value = 406
$ws.Range("K19").Value = -29.556650246305
$ws.Range("L19").Value = -37.826086956521
$ws.Range("M19").Value = -6.229508196721
$ws.Range("N19").Value = -24.538258575197

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 83
$ws.Range("J20").Value = 127
$ws.Range("K20").Value = -34.645669291338
$ws.Range("L20").Value = -45.394736842105
$ws.Range("M20").Value = -17
$ws.Range("N20").Value = -89.890377588306

# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -37.5
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = -28.205128205128
$ws.Range("I21").Value = 811
$ws.Range("J21").Value = 1134
$ws.Range("K21").Value = -28.483245149911
$ws.Range("L21").Value = -26.473254759746
$ws.Range("M21").Value = -18.737474949899
$ws.Range("N21").Value = -76.668584579977

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -38.461538461538
$ws.Range("F24").Value = 52
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = -32.467532467532
$ws.Range("I24").Value = 594
$ws.Range("J24").Value = 795
$ws.Range("K24").Value = -25.283018867924
$ws.Range("L24").Value = -42.994241842610
$ws.Range("M24").Value = -8.755760368663

# Row 25
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -46.153846153846
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -60
$ws.Range("I25").Value = 179
$ws.Range("J25").Value = 323
$ws.Range("K25").Value = -44.582043343653
$ws.Range("L25").Value = -69.661016949152

# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -15.384615384615
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 52
$ws.Range("H26").Value = 3.846153846153
$ws.Range("I26").Value = 441
$ws.Range("J26").Value = 461
$ws.Range("K26").Value = -4.338394793926
$ws.Range("L26").Value = 2.083333333333
$ws.Range("M26").Value = -16.949152542372

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 27
$ws.Range("K27").Value = 8
$ws.Range("L27").Value = 28.571428571428

# Row 28
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 4
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 48
$ws.Range("J28").Value = 48
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -30.434782608695

# Row 29
$ws.Range("M29").Value = -62.5
$ws.Range("N29").Value = -88

# Row 30
$ws.Range("M30").Value = -57.142857142857
$ws.Range("N30").Value = -87.5

